$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 19 - "Git": merge the two runs "durch " and "Linus " (identical
# formatting) into a single run "durch Linus ".
# ---------------------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$shape19 = $s19.Shapes.Item(2)
$tr19 = $shape19.TextFrame.TextRange
$para19 = $tr19.Paragraphs(2, 1)
$start19 = $para19.Start
$paraText19 = $para19.Text
$offLinus = $paraText19.IndexOf("Linus ")
$offDurch = $paraText19.IndexOf("durch ")

# Delete the "Linus " run entirely...
$linusRange = $tr19.Characters($start19 + $offLinus, 6)
$linusRange.Text = ""

# ...then re-insert it right after "durch " so it gets folded into that run.
$durchRange = $tr19.Characters($start19 + $offDurch, 6)
$durchRange.InsertAfter("Linus ") | Out-Null

# ---------------------------------------------------------------------------
# Slide 4 - "CI - Grundsätze": split the single run describing the
# pre-integration tests so "Statischen" becomes "statische" as its own run.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shape4 = $s4.Shapes.Item(2)
$tr4 = $shape4.TextFrame.TextRange
$para4 = $tr4.Paragraphs(4, 1)
$start4 = $para4.Start
$paraText4 = $para4.Text
$offStat = $paraText4.IndexOf("Statischen ")
$statRange = $tr4.Characters($start4 + $offStat, 11)
$statRange.Text = "statische "

# ---------------------------------------------------------------------------
# Slide 5 - "CI - Grundsätze": reword "Bspw. 1 Mal täglich mindestens" into
# "Bspw. mindestens 1 Mal täglich".
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shape5 = $s5.Shapes.Item(2)
$tr5 = $shape5.TextFrame.TextRange
$para5 = $tr5.Paragraphs(5, 1)
$start5 = $para5.Start

# Drop the trailing " mindestens" first so earlier offsets stay valid.
$tailRange = $tr5.Characters($start5 + 19, 11)
$tailRange.Text = ""

# Replace "1 " with "mindestens 1 " so the word moves in front of "Mal".
$oneRange = $tr5.Characters($start5 + 6, 2)
$oneRange.Text = "mindestens 1 "

# Finally, split "täglich" off into its own run.
$para5b = $tr5.Paragraphs(5, 1)
$paraText5b = $para5b.Text
$offTaeglich = $paraText5b.IndexOf("täglich")
$taeglichRange = $tr5.Characters($start5 + $offTaeglich, 7)
$taeglichRange.Text = "täglich"
